$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap header text between K1 and L1 (Postcode <-> Street Name) ---
$k1 = $ws.Range("K1").Value()
$l1 = $ws.Range("L1").Value()
$ws.Range("K1").Value = $l1
$ws.Range("L1").Value = $k1

# --- Swap the column widths of columns K (11) and L (12) ---
$ws.Columns.Item(12).ColumnWidth = 13.166666666666666
$ws.Columns.Item(11).ColumnWidth = 14.5

# --- Fill in new row 6 test data ---
$ws.Range("D6").Value = 45171943
$ws.Range("H6").Value = "ADSSN"
$ws.Range("I6").Value = "ADSFN"
$ws.Range("K6").Value = "PARRAMATTA"
$ws.Range("J6").Value = "'10/10/1975"
$ws.Range("J6").NumberFormat = "mm-dd-yy"
$ws.Range("L6").Value = 2140

# --- Update the active cell selection ---
$ws.Range("J6").Select()
